# Refresh "latest" optimisation output (run 70): the rolling forecast
# horizon advances by one 30-minute period, so the oldest "Detailed" row
# (the prior historical interval) drops off, every remaining row's
# DateTime/Pump_Status shift up one slot, a freshly re-forecast Price/Type
# series is written in, and the dependent Schedule cost figures are
# recomputed for the new horizon.

$wb = $excel.ActiveWorkbook
$schedule = $wb.Worksheets.Item("Schedule")
$detailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet: update Cost ($) and Unit Cost ($/ML) for rows 2-3 ---
$schedule.Range("E2").Value = -56.21536049999997
$schedule.Range("F2").Value = -1.239315707671957
$schedule.Range("E3").Value = 371.6685765
$schedule.Range("F3").Value = 24.58125505952381

# --- Detailed sheet: drop oldest timestamp row, shift data up one period ---
$detailed.Rows("2:2").Delete()

# --- Detailed sheet: write refreshed Price/Type values (and re-affirm DateTime/Pump_Status) for rows 2-48 ---
$detailed.Range("A2").Value = 46040.02083333334
$detailed.Range("B2").Value = 36.06032
$detailed.Range("C2").Value = "historical"
$detailed.Range("E2").Value = "OFF"
$detailed.Range("A3").Value = 46040.04166666666
$detailed.Range("B3").Value = 35.88
$detailed.Range("C3").Value = "historical"
$detailed.Range("E3").Value = "OFF"
$detailed.Range("A4").Value = 46040.0625
$detailed.Range("B4").Value = 36.2
$detailed.Range("C4").Value = "forecast"
$detailed.Range("E4").Value = "OFF"
$detailed.Range("A5").Value = 46040.08333333334
$detailed.Range("B5").Value = 57.06019
$detailed.Range("C5").Value = "forecast"
$detailed.Range("E5").Value = "OFF"
$detailed.Range("A6").Value = 46040.10416666666
$detailed.Range("B6").Value = 57.0602
$detailed.Range("C6").Value = "forecast"
$detailed.Range("E6").Value = "OFF"
$detailed.Range("A7").Value = 46040.125
$detailed.Range("B7").Value = 57.06021
$detailed.Range("C7").Value = "forecast"
$detailed.Range("E7").Value = "OFF"
$detailed.Range("A8").Value = 46040.14583333334
$detailed.Range("B8").Value = 57.06022
$detailed.Range("C8").Value = "forecast"
$detailed.Range("E8").Value = "OFF"
$detailed.Range("A9").Value = 46040.16666666666
$detailed.Range("B9").Value = 57.06022
$detailed.Range("C9").Value = "forecast"
$detailed.Range("E9").Value = "OFF"
$detailed.Range("A10").Value = 46040.1875
$detailed.Range("B10").Value = 57.06003
$detailed.Range("C10").Value = "forecast"
$detailed.Range("E10").Value = "OFF"
$detailed.Range("A11").Value = 46040.20833333334
$detailed.Range("B11").Value = 57.06003
$detailed.Range("C11").Value = "forecast"
$detailed.Range("E11").Value = "OFF"
$detailed.Range("A12").Value = 46040.22916666666
$detailed.Range("B12").Value = 57.06003
$detailed.Range("C12").Value = "forecast"
$detailed.Range("E12").Value = "OFF"
$detailed.Range("A13").Value = 46040.25
$detailed.Range("B13").Value = 57.06003
$detailed.Range("C13").Value = "forecast"
$detailed.Range("E13").Value = "OFF"
$detailed.Range("A14").Value = 46040.27083333334
$detailed.Range("B14").Value = 36.06
$detailed.Range("C14").Value = "forecast"
$detailed.Range("E14").Value = "OFF"
$detailed.Range("A15").Value = 46040.29166666666
$detailed.Range("B15").Value = 36.0595
$detailed.Range("C15").Value = "forecast"
$detailed.Range("E15").Value = "ON"
$detailed.Range("A16").Value = 46040.3125
$detailed.Range("B16").Value = 11.56192
$detailed.Range("C16").Value = "forecast"
$detailed.Range("E16").Value = "ON"
$detailed.Range("A17").Value = 46040.33333333334
$detailed.Range("B17").Value = -5.01
$detailed.Range("C17").Value = "forecast"
$detailed.Range("E17").Value = "ON"
$detailed.Range("A18").Value = 46040.35416666666
$detailed.Range("B18").Value = 0.5954700000000001
$detailed.Range("C18").Value = "forecast"
$detailed.Range("E18").Value = "ON"
$detailed.Range("A19").Value = 46040.375
$detailed.Range("B19").Value = 0.009549999999999999
$detailed.Range("C19").Value = "forecast"
$detailed.Range("E19").Value = "ON"
$detailed.Range("A20").Value = 46040.39583333334
$detailed.Range("B20").Value = -3.6481
$detailed.Range("C20").Value = "forecast"
$detailed.Range("E20").Value = "ON"
$detailed.Range("A21").Value = 46040.41666666666
$detailed.Range("B21").Value = -4.55365
$detailed.Range("C21").Value = "forecast"
$detailed.Range("E21").Value = "ON"
$detailed.Range("A22").Value = 46040.4375
$detailed.Range("B22").Value = 0
$detailed.Range("C22").Value = "forecast"
$detailed.Range("E22").Value = "ON"
$detailed.Range("A23").Value = 46040.45833333334
$detailed.Range("B23").Value = -5.51
$detailed.Range("C23").Value = "forecast"
$detailed.Range("E23").Value = "ON"
$detailed.Range("A24").Value = 46040.47916666666
$detailed.Range("B24").Value = -1.31495
$detailed.Range("C24").Value = "forecast"
$detailed.Range("E24").Value = "ON"
$detailed.Range("A25").Value = 46040.5
$detailed.Range("B25").Value = -4.50424
$detailed.Range("C25").Value = "forecast"
$detailed.Range("E25").Value = "ON"
$detailed.Range("A26").Value = 46040.52083333334
$detailed.Range("B26").Value = 0
$detailed.Range("C26").Value = "forecast"
$detailed.Range("E26").Value = "ON"
$detailed.Range("A27").Value = 46040.54166666666
$detailed.Range("B27").Value = -5.51011
$detailed.Range("C27").Value = "forecast"
$detailed.Range("E27").Value = "ON"
$detailed.Range("A28").Value = 46040.5625
$detailed.Range("B28").Value = -6.8
$detailed.Range("C28").Value = "forecast"
$detailed.Range("E28").Value = "ON"
$detailed.Range("A29").Value = 46040.58333333334
$detailed.Range("B29").Value = -14
$detailed.Range("C29").Value = "forecast"
$detailed.Range("E29").Value = "ON"
$detailed.Range("A30").Value = 46040.60416666666
$detailed.Range("B30").Value = -23.5
$detailed.Range("C30").Value = "forecast"
$detailed.Range("E30").Value = "ON"
$detailed.Range("A31").Value = 46040.625
$detailed.Range("B31").Value = -14
$detailed.Range("C31").Value = "forecast"
$detailed.Range("E31").Value = "ON"
$detailed.Range("A32").Value = 46040.64583333334
$detailed.Range("B32").Value = -14
$detailed.Range("C32").Value = "forecast"
$detailed.Range("E32").Value = "ON"
$detailed.Range("A33").Value = 46040.66666666666
$detailed.Range("B33").Value = -7.1156
$detailed.Range("C33").Value = "forecast"
$detailed.Range("E33").Value = "ON"
$detailed.Range("A34").Value = 46040.6875
$detailed.Range("B34").Value = -6.51112
$detailed.Range("C34").Value = "forecast"
$detailed.Range("E34").Value = "ON"
$detailed.Range("A35").Value = 46040.70833333334
$detailed.Range("B35").Value = -6
$detailed.Range("C35").Value = "forecast"
$detailed.Range("E35").Value = "ON"
$detailed.Range("A36").Value = 46040.72916666666
$detailed.Range("B36").Value = 0.63476
$detailed.Range("C36").Value = "forecast"
$detailed.Range("E36").Value = "ON"
$detailed.Range("A37").Value = 46040.75
$detailed.Range("B37").Value = 3.95694
$detailed.Range("C37").Value = "forecast"
$detailed.Range("E37").Value = "ON"
$detailed.Range("A38").Value = 46040.77083333334
$detailed.Range("B38").Value = 11.50285
$detailed.Range("C38").Value = "forecast"
$detailed.Range("E38").Value = "ON"
$detailed.Range("A39").Value = 46040.79166666666
$detailed.Range("B39").Value = 36.2
$detailed.Range("C39").Value = "forecast"
$detailed.Range("E39").Value = "OFF"
$detailed.Range("A40").Value = 46040.8125
$detailed.Range("B40").Value = 53.65862
$detailed.Range("C40").Value = "forecast"
$detailed.Range("E40").Value = "OFF"
$detailed.Range("A41").Value = 46040.83333333334
$detailed.Range("B41").Value = 56.98
$detailed.Range("C41").Value = "forecast"
$detailed.Range("E41").Value = "ON"
$detailed.Range("A42").Value = 46040.85416666666
$detailed.Range("B42").Value = 56.98
$detailed.Range("C42").Value = "forecast"
$detailed.Range("E42").Value = "ON"
$detailed.Range("A43").Value = 46040.875
$detailed.Range("B43").Value = 56.98
$detailed.Range("C43").Value = "forecast"
$detailed.Range("E43").Value = "ON"
$detailed.Range("A44").Value = 46040.89583333334
$detailed.Range("B44").Value = 45.73237
$detailed.Range("C44").Value = "forecast"
$detailed.Range("E44").Value = "ON"
$detailed.Range("A45").Value = 46040.91666666666
$detailed.Range("B45").Value = 56.34597
$detailed.Range("C45").Value = "forecast"
$detailed.Range("E45").Value = "ON"
$detailed.Range("A46").Value = 46040.9375
$detailed.Range("B46").Value = 36.0601
$detailed.Range("C46").Value = "forecast"
$detailed.Range("E46").Value = "ON"
$detailed.Range("A47").Value = 46040.95833333334
$detailed.Range("B47").Value = 36.0601
$detailed.Range("C47").Value = "forecast"
$detailed.Range("E47").Value = "ON"
$detailed.Range("A48").Value = 46040.97916666666
$detailed.Range("B48").Value = 36.06
$detailed.Range("C48").Value = "forecast"
$detailed.Range("E48").Value = "ON"
